$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.66539192199707
$ws.Range("B1").Value = 3.821148633956909
$ws.Range("C1").Value = 2.069429874420166
$ws.Range("D1").Value = 1.368826389312744
$ws.Range("E1").Value = 1.128280401229858
